# Planificacion.xlsx - "Enemigos" sheet update
# Commit: Agregado spawneo de 3 tipos distintos de enemigos en clusters
#
# The columns E,F,G (header row 2 + data rows 3-5) get rotated:
#   new E <- old F   (Daño)
#   new F <- old G   (Vida)
#   new G <- old E   (Rango Seguimiento)
# Row 5's new E value is bumped from 3 -> 3.5 (a genuine data change on
# top of the rotation).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- capture current (pre-edit) header + data values for E:G, rows 2-5 ---
$oldE2 = $ws.Range("E2").Value2
$oldF2 = $ws.Range("F2").Value2
$oldG2 = $ws.Range("G2").Value2

$oldE3 = $ws.Range("E3").Value2
$oldF3 = $ws.Range("F3").Value2
$oldG3 = $ws.Range("G3").Value2

$oldE4 = $ws.Range("E4").Value2
$oldF4 = $ws.Range("F4").Value2
$oldG4 = $ws.Range("G4").Value2

$oldE5 = $ws.Range("E5").Value2
$oldF5 = $ws.Range("F5").Value2
$oldG5 = $ws.Range("G5").Value2

# --- write back the rotated values ---
$ws.Range("E2").Value = $oldF2
$ws.Range("F2").Value = $oldG2
$ws.Range("G2").Value = $oldE2

$ws.Range("E3").Value = $oldF3
$ws.Range("F3").Value = $oldG3
$ws.Range("G3").Value = $oldE3

$ws.Range("E4").Value = $oldF4
$ws.Range("F4").Value = $oldG4
$ws.Range("G4").Value = $oldE4

# Row 5: same rotation, but the value landing in E5 is bumped to 3.5
$ws.Range("E5").Value = 3.5
$ws.Range("F5").Value = $oldG5
$ws.Range("G5").Value = $oldE5

# --- match the resulting bestFit column widths (they follow the moved content) ---
$ws.Columns.Item(5).ColumnWidth = 5.5703125
$ws.Columns.Item(6).ColumnWidth = 5
$ws.Columns.Item(7).ColumnWidth = 18.28515625

# --- update the active selection to H11 (as recorded in the saved view state) ---
$ws.Range("H11").Select()
